# Fix missing expert testimony achievement in comprehensive resumes:
# insert four new bullet paragraphs into the "KEY ACHIEVEMENTS AND IMPACT"
# section, right after the "Platform impact" bullet and before the
# "TECHNICAL SKILLS" heading.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$bullet = [char]0x2022
$plusMinus = [char]0xB1

$newParagraphs = @(
    "$bullet Real-time collaboration at national scale",
    "$bullet Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ${plusMinus}4.2% to ${plusMinus}2.1%",
    "$bullet Increased voter turnout prediction accuracy from 71% to 87%",
    "$bullet Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
)

$insertPos = $anchor.End
foreach ($text in $newParagraphs) {
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertParagraphAfter()
    $newTextRange = $d.Range($insertPos + 1, $insertPos + 1)
    $newTextRange.Text = $text
    $insertPos = $insertPos + 1 + $text.Length
}
